$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.822.92"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.637.06"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.15%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.55"

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5057"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.24%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  +0.37%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06421"
$ws.Range("E9").Value = "  +1.30%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  +0.36%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").Value = "  +0.70%  "

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.279"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.641.34"
$ws.Range("E13").Value = "  +0.51%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.863.31"
$ws.Range("E14").Value = "  +0.21%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5615"
$ws.Range("E15").Value = "  +3.45%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0₅7590"
$ws.Range("E16").Value = "  -1.50%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.11"
$ws.Range("E17").Value = "  -1.37%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.840.26"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.13%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.53"
$ws.Range("E20").Value = "  -0.17%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.313"
$ws.Range("E21").Value = "  -2.38%  "

# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.875"
$ws.Range("E22").Value = "  -0.19%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.117"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24 - BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.771"
$ws.Range("E25").Value = "  -6.27%  "

# Rows 26/27 - Stellar and Monero swap ranking positions
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.09"
$ws.Range("E26").Value = "  -1.85%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1266"
$ws.Range("E27").Value = "  +1.66%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.796"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.44"
$ws.Range("E29").Value = "  -0.81%  "

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("E30").Value = "  +0.66%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04866"
$ws.Range("E31").Value = "  +0.46%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.297"
$ws.Range("E32").Value = "  +2.07%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.217"
$ws.Range("E33").Value = "  +1.27%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.47%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.376"
$ws.Range("E35").Value = "  +0.15%  "

# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9018"
$ws.Range("E36").Value = "  -0.39%  "

# Row 37 - MXToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.582"
$ws.Range("E37").Value = "  +0.32%  "

# Row 38 - ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5514"
$ws.Range("E38").Value = "  +0.77%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.125.64"
$ws.Range("E39").Value = "  -0.01%  "

# Row 40 - VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01559"
$ws.Range("E40").Value = "  +0.42%  "

# Row 41 - PaxDollar
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9961"
$ws.Range("E41").Value = "  -0.51%  "

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.532"
$ws.Range("E42").Value = "  -0.72%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8010"
$ws.Range("E43").Value = "  +0.08%  "

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.81"
$ws.Range("E44").Value = "  -0.54%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.775.80"
$ws.Range("E45").Value = "  +0.29%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -7.97%  "

# Row 47 - Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.40"
$ws.Range("E47").Value = "  +1.03%  "

# Row 48 - Mantle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4294"
$ws.Range("E48").Value = "  -4.03%  "

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.720"
$ws.Range("E49").Value = "  +2.46%  "

# Row 50 - Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05048"
$ws.Range("E50").Value = "  -2.19%  "

# Row 51 - Frax
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.18%  "
